# Add a "time_taken" column (F) to the panel data sheet, mirroring the
# header style used by the other header cells and filling each data row
# with its recorded timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy the format from the existing header (E1) so it
# matches the bold/centered look of the other headers, then set its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps for F2:F13
$timestamps = @(
    "2021-10-05 10:50:53.219482",
    "2021-10-05 10:50:53.219493",
    "2021-10-05 10:50:53.219497",
    "2021-10-05 10:50:53.219499",
    "2021-10-05 10:50:53.219502",
    "2021-10-05 10:50:53.219505",
    "2021-10-05 10:50:53.219507",
    "2021-10-05 10:50:53.219510",
    "2021-10-05 10:50:53.219513",
    "2021-10-05 10:50:53.219515",
    "2021-10-05 10:50:53.219518",
    "2021-10-05 10:50:53.219520"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
